$wb = $excel.ActiveWorkbook
$choices = $wb.Worksheets.Item("choices")

$choices.Range("B2").Value = "a1"
$choices.Range("B3").Value = "a0"
$choices.Range("B4").Value = "a1"
$choices.Range("B5").Value = "a0"
$choices.Range("B6").Value = "a888"
$choices.Range("B7").Value = "a888"
$choices.Range("B8").Value = "a1"
$choices.Range("B9").Value = "a2"
$choices.Range("B10").Value = "a3"
$choices.Range("B11").Value = "a4"
$choices.Range("B12").Value = "a5"
$choices.Range("B13").Value = "a6"
$choices.Range("B14").Value = "a7"
$choices.Range("B15").Value = "a888"
$choices.Range("B16").Value = "a1"
$choices.Range("B17").Value = "a2"
$choices.Range("B18").Value = "a3"
$choices.Range("B19").Value = "a4"
$choices.Range("B20").Value = "a5"
$choices.Range("B21").Value = "a6"
$choices.Range("B22").Value = "a7"
$choices.Range("B23").Value = "a8"
$choices.Range("B24").Value = "a9"
$choices.Range("B25").Value = "a10"
$choices.Range("B26").Value = "a888"
$choices.Range("B27").Value = "a0"
$choices.Range("B28").Value = "a1"
$choices.Range("B29").Value = "a2"
$choices.Range("B30").Value = "a3"
$choices.Range("B31").Value = "a4"
$choices.Range("B32").Value = "a5"
$choices.Range("B33").Value = "a6"
$choices.Range("B34").Value = "a7"
$choices.Range("B35").Value = "a8"
$choices.Range("B36").Value = "a9"
$choices.Range("B37").Value = "a888"
$choices.Range("B38").Value = "a1"
$choices.Range("B39").Value = "a2"
$choices.Range("B40").Value = "a3"
$choices.Range("B41").Value = "a888"
$choices.Range("B42").Value = "a0"
$choices.Range("B43").Value = "a1"
$choices.Range("B44").Value = "a888"
$choices.Range("B45").Value = "a999"
$choices.Range("B46").Value = "a1"
$choices.Range("B47").Value = "a0"
$choices.Range("B48").Value = "a888"
$choices.Range("B49").Value = "a999"
$choices.Range("B50").Value = "a1"
$choices.Range("B51").Value = "a2"
$choices.Range("B52").Value = "a3"
$choices.Range("B53").Value = "a4"
$choices.Range("B59").Value = "a0"
$choices.Range("B60").Value = "a1"
$choices.Range("B61").Value = "a2"
$choices.Range("B62").Value = "a3"
$choices.Range("B63").Value = "a4"
$choices.Range("B64").Value = "a5"
$choices.Range("B65").Value = "a888"
$choices.Range("B66").Value = "a9999"
$choices.Range("B67").Value = "a0"
$choices.Range("B68").Value = "a1"
$choices.Range("B69").Value = "a2"
$choices.Range("B70").Value = "a3"
$choices.Range("B71").Value = "a4"
$choices.Range("B72").Value = "a5"
$choices.Range("B73").Value = "a6"
$choices.Range("B74").Value = "a7"
$choices.Range("B75").Value = "a8"
$choices.Range("B76").Value = "a9"
$choices.Range("B77").Value = "a888"
$choices.Range("B78").Value = "a1"
$choices.Range("B79").Value = "a1"
$choices.Range("B80").Value = "a2"
$choices.Range("B81").Value = "a888"
$choices.Range("B82").Value = "a1"
$choices.Range("B83").Value = "a0"

$choices.Range("B14:B15").HorizontalAlignment = 1

$choices.Activate()
$choices.Range("B6").Select()
